$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 8 (year 2025) metrics per commit "atualizei dados bibi e add"
$ws.Range("C8").Value = 1274
$ws.Range("D8").Value = 206
$ws.Range("E8").Value = 1068
$ws.Range("F8").Value = 8.44954881050041
$ws.Range("G8").Value = 83.83045525902669
$ws.Range("H8").Value = 16.16954474097331
